$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 27.80961241492737
$ws.Range("B3").Value = 25.42882193453775
$ws.Range("B4").Value = 20.14878607324049
$ws.Range("B5").Value = 9.653878052057268
$ws.Range("B6").Value = 7.394374712868707
$ws.Range("B7").Value = 5.019209798603036
$ws.Range("B8").Value = 4.545317013765384
